$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.957.61'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').Value = '3.163.11'
$ws.Range('E3').Value = '  -7.59%  '
$ws.Range('D5').Value = '568.96'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').Value = '170.01'
$ws.Range('E6').Value = '  -6.33%  '
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '3.161.63'
$ws.Range('E9').Value = '  -7.65%  '
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  -5.70%  '
$ws.Range('E11').Value = '  -6.01%  '
$ws.Range('E12').Value = '  -4.74%  '
$ws.Range('D13').Value = '3.709.13'
$ws.Range('E13').Value = '  -7.77%  '
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '27.15'
$ws.Range('E15').Value = '  -7.16%  '
$ws.Range('D16').Value = '64.893.19'
$ws.Range('E16').Value = '  -2.54%  '
$ws.Range('E17').Value = '  -6.22%  '
$ws.Range('D18').Value = '3.163.69'
$ws.Range('E18').Value = '  -7.68%  '
$ws.Range('D19').Value = '5.75'
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('D20').Value = '12.83'
$ws.Range('E20').Value = '  -7.34%  '
$ws.Range('D21').Value = '357.68'
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').Value = '7.29'
$ws.Range('E22').Value = '  -4.19%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '69.27'
$ws.Range('E24').Value = '  -5.36%  '
$ws.Range('E25').Value = '  -6.85%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.304.17'
$ws.Range('E26').Value = '  -7.71%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0000116'
$ws.Range('E27').Value = '  -7.51%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '9.75'
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.176'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.91'
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '22.08'
$ws.Range('E33').Value = '  -5.42%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '5.33'
$ws.Range('E34').Value = '  -8.08%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '1.21'
$ws.Range('E35').Value = '  -4.57%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '6.64'
$ws.Range('E36').Value = '  -5.97%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.45'
$ws.Range('E37').Value = '  -6.65%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '155.84'
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  -3.69%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.78'
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '26.13'
$ws.Range('E41').Value = '  -5.33%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.666.68'
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '2.47'
$ws.Range('E43').Value = '  -6.80%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '4.19'
$ws.Range('E44').Value = '  -5.14%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '6.04'
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '39.47'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0658'
$ws.Range('E47').Value = '  -4.43%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '24.06'
$ws.Range('E48').Value = '  -3.84%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '325.03'
$ws.Range('E49').Value = '  -3.90%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0273'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.103'
$ws.Range('E51').Value = '  -1.51%  '
